{"js": "// The commit corrects a typo'd student ID in the authors line:\n//   \"Danilo Verde 1000001353, Nicol\u00f2 Mazzola\"\n//     -> \"Danilo Verde 1000069238, Nicol\u00f2 Mazzola\"\n// (the trailing \" 1000067652\" that belongs to the second author is a\n// separate run and must stay untouched).\n//\n// Locate the exact old substring and replace it with the corrected one so\n// the edit is scoped tightly to the changed digits, leaving every other\n// run/paragraph in the document exactly as it was.\n\nconst body = context.document.body;\n\nconst searchResults = body.search(\"1000001353\", { matchCase: true });\nsearchResults.load(\"items\");\nawait context.sync();\n\nif (searchResults.items.length > 0) {\n  searchResults.items[0].insertText(\"1000069238\", \"Replace\");\n} else {\n  // Fallback in case the surrounding text already changed shape: replace\n  // just the differing middle digits so \"Danilo Verde 10000\" + \"69238\" +\n  // \", Nicol\u00f2 Mazzola\" is preserved either way.\n  const narrow = body.search(\"01353\", { matchCase: true });\n  narrow.load(\"items\");\n  await context.sync();\n  if (narrow.items.length > 0) {\n    narrow.items[0].insertText(\"69238\", \"Replace\");\n  }\n}\n\nawait context.sync();\n", "ps1": "# The commit corrects a typo'd student ID in the authors line:\n#   \"Danilo Verde 1000001353, Nicol\u00f2 Mazzola\"\n#     -> \"Danilo Verde 1000069238, Nicol\u00f2 Mazzola\"\n# (the trailing \" 1000067652\" that belongs to the second author sits in its\n# own run and must be left untouched).\n#\n# Use Find/Replace scoped to the exact old ID so only those digits change;\n# everything else in the document is left byte-for-byte identical.\n\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.Text = \"1000001353\"\n$find.Replacement.Text = \"1000069238\"\n$found = $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n\nif (-not $found) {\n    # Fallback: only the differing middle digits, in case the surrounding\n    # text already changed shape.\n    $find2 = $d.Content.Find\n    $find2.Text = \"01353\"\n    $find2.Replacement.Text = \"69238\"\n    $find2.Execute($find2.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find2.Replacement.Text, 2)\n}\n"}
